$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    "greifen",
    "schneiden",
    "zögern",
    "liefern",
    "schreiten",
    "gelten",
    "trennen",
    "schrecken",
    "plaudern",
    "spinnen",
    "knurren",
    "runden",
    "schwingen",
    "stillen",
    "altern",
    "tropfen",
    "freuen",
    "rauben",
    "filmen",
    "stärken",
    "zeigen",
    "wehtun",
    "weichen",
    "testen",
    "brauen",
    "laufen",
    "löschen",
    "bauen",
    "wüten",
    "machen",
    "wahren",
    "wundern",
    "bluten",
    "führen",
    "malen",
    "stehlen",
    "sterben",
    "binden",
    "helfen",
    "kratzen",
    "segeln",
    "pflanzen",
    "foltern",
    "heulen",
    "scheinen",
    "schämen",
    "nerven",
    "wenden",
    "loben",
    "klingeln",
    "tollen",
    "locken",
    "gnaden",
    "schwimmen",
    "starten",
    "kümmern",
    "kichern",
    "drehen",
    "füllen",
    "fischen",
    "jagen",
    "sinken",
    "münzen",
    "hupen",
    "geben",
    "mauern",
    "ärgern",
    "suchen",
    "sichern",
    "lehnen",
    "klingen",
    "knarren",
    "lesen",
    "klettern",
    "ziehen",
    "mögen",
    "leugnen",
    "treiben",
    "hauen",
    "rufen",
    "brauchen",
    "heben",
    "schicken",
    "dringen",
    "warnen",
    "zünden",
    "pissen",
    "ändern",
    "bitten",
    "planen",
    "stammen",
    "danken",
    "leeren",
    "achten",
    "rechnen",
    "spielen",
    "regeln",
    "hören",
    "klagen",
    "graben",
    "räumen",
    "fallen",
    "wachsen",
    "landen",
    "rühren",
    "melden",
    "weinen",
    "ehren",
    "rasen",
    "schulden",
    "wollen",
    "wirken",
    "streichen",
    "backen",
    "saufen",
    "schnellen",
    "schmecken",
    "schenken",
    "gründen",
    "trotzen",
    "drücken",
    "bergen",
    "buchen",
    "kehren",
    "hassen",
    "fließen",
    "bellen",
    "boxen",
    "heilen",
    "reizen",
    "lügen",
    "spannen",
    "lockern",
    "feuern",
    "rücken",
    "stimmen",
    "feiern",
    "fällen",
    "irren",
    "leiden",
    "scheitern",
    "kosten",
    "spüren",
    "schütteln",
    "grüßen",
    "arten",
    "quälen",
    "jubeln",
    "meinen",
    "töten",
    "tauchen",
    "flüchten",
    "sprengen",
    "opfern",
    "fangen",
    "deuten",
    "liegen",
    "sperren",
    "schwören",
    "dienen",
    "schauen",
    "pfeifen",
    "fahren",
    "betteln",
    "schlucken",
    "stürmen",
    "platzen",
    "flehen",
    "trauen",
    "wetten",
    "biegen",
    "werden",
    "steuern",
    "ächzen",
    "zielen",
    "enden",
    "formen",
    "erben",
    "duschen",
    "folgen",
    "äußern",
    "schwächen",
    "decken",
    "sammeln",
    "teilen",
    "passen",
    "sorgen",
    "stecken",
    "werfen",
    "seufzen",
    "schultern",
    "siegen"
)

for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newWords[$i]
}